$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{A="ECs"; B="Anxa1"; C="Fpr2"; D="FAPs"; E=3; F=1; G=28.49436566666666; H=85.48309699999999; I=0.04826188997421807; J=0.04826188997421808; K=2; L=0.6666666666666666; M=1.109174333333333; N=3.327523; O=0.06061832081580493; P=0.06061832081580493; Q=31.60521904208122; R=284.446971378731; S=0.00292555472963423; T=0.002925554729634231},
  @{A="ECs"; B="Anxa1"; C="Fpr2"; D="M1"; E=3; F=1; G=28.49436566666666; H=85.48309699999999; I=0.04826188997421807; J=0.04826188997421808; K=3; L=1; M=4.402069333333333; N=13.206208; O=0.2405808023879173; P=0.2405808023879173; Q=125.4341732740195; R=1128.907559466176; S=0.01161088421475476; T=0.01161088421475476},
  @{A="ECs"; B="Anxa1"; C="Fpr2"; D="M2"; E=3; F=1; G=28.49436566666666; H=85.48309699999999; I=0.04826188997421807; J=0.04826188997421808; K=3; L=1; M=12.78643133333333; N=38.359294; O=0.6988008767962779; P=0.6988008767962779; Q=364.3412499837241; R=3279.071249853517; S=0.03372545102982908; T=0.03372545102982909},
  @{A="FAPs"; B="Anxa1"; C="Fpr2"; D="FAPs"; E=3; F=1; G=142.0814363333334; H=426.244309; I=0.2406482294750577; J=0.2406482294750577; K=2; L=0.6666666666666666; M=1.109174333333333; N=3.327523; O=0.06061832081580493; P=0.06061832081580493; Q=157.5930824240675; R=1418.337741816607; S=0.01458769157807449; T=0.01458769157807449},
  @{A="FAPs"; B="Anxa1"; C="Fpr2"; D="M1"; E=3; F=1; G=142.0814363333334; H=426.244309; I=0.2406482294750577; J=0.2406482294750577; K=3; L=1; M=4.402069333333333; N=13.206208; O=0.2405808023879173; P=0.2405808023879173; Q=625.4523337189192; R=5629.071003470272; S=0.05789534414034102; T=0.05789534414034102},
  @{A="FAPs"; B="Anxa1"; C="Fpr2"; D="M2"; E=3; F=1; G=142.0814363333334; H=426.244309; I=0.2406482294750577; J=0.2406482294750577; K=3; L=1; M=12.78643133333333; N=38.359294; O=0.6988008767962779; P=0.6988008767962779; Q=1816.714529417539; R=16350.43076475785; S=0.1681651937566422; T=0.1681651937566422},
  @{A="M1"; B="Anxa1"; C="Fpr2"; D="FAPs"; E=3; F=1; G=209.9177196666667; H=629.753159; I=0.3555448824061003; J=0.3555448824061003; K=2; L=0.6666666666666666; M=1.109174333333333; N=3.327523; O=0.06061832081580493; P=0.06061832081580493; Q=232.8353467661285; R=2095.518120895157; S=0.02155253374611062; T=0.02155253374611062},
  @{A="M1"; B="Anxa1"; C="Fpr2"; D="M1"; E=3; F=1; G=209.9177196666667; H=629.753159; I=0.3555448824061003; J=0.3555448824061003; K=3; L=1; M=4.402069333333333; N=13.206208; O=0.2405808023879173; P=0.2405808023879173; Q=924.0723562678969; R=8316.651206411072; S=0.08553727309417729; T=0.08553727309417729},
  @{A="M1"; B="Anxa1"; C="Fpr2"; D="M2"; E=3; F=1; G=209.9177196666667; H=629.753159; I=0.3555448824061003; J=0.3555448824061003; K=3; L=1; M=12.78643133333333; N=38.359294; O=0.6988008767962779; P=0.6988008767962779; Q=2684.09850816775; R=24156.88657350975; S=0.2484550755658124; T=0.2484550755658124},
  @{A="M2"; B="Anxa1"; C="Fpr2"; D="FAPs"; E=3; F=1; G=132.271009; H=396.813027; I=0.2240319703135046; J=0.2240319703135046; K=2; L=0.6666666666666666; M=1.109174333333333; N=3.327523; O=0.06061832081580493; P=0.06061832081580493; Q=146.7116082269023; R=1320.404474042121; S=0.0135804418494609; T=0.0135804418494609},
  @{A="M2"; B="Anxa1"; C="Fpr2"; D="M1"; E=3; F=1; G=132.271009; H=396.813027; I=0.2240319703135046; J=0.2240319703135046; K=3; L=1; M=4.402069333333333; N=13.206208; O=0.2405808023879173; P=0.2405808023879173; Q=582.2661524079574; R=5240.395371671617; S=0.05389779117856899; T=0.05389779117856899},
  @{A="M2"; B="Anxa1"; C="Fpr2"; D="M2"; E=3; F=1; G=132.271009; H=396.813027; I=0.2240319703135046; J=0.2240319703135046; K=3; L=1; M=12.78643133333333; N=38.359294; O=0.6988008767962779; P=0.6988008767962779; Q=1691.274173969216; R=15221.46756572294; S=0.1565537372854747; T=0.1565537372854747},
  @{A="sCs"; B="Anxa1"; C="Fpr2"; D="FAPs"; E=3; F=1; G=77.646779; H=232.940337; I=0.1315130278311194; J=0.1315130278311194; K=2; L=0.6666666666666666; M=1.109174333333333; N=3.327523; O=0.06061832081580493; P=0.06061832081580493; Q=86.12381433280565; R=775.114328995251; S=0.007972098912524678; T=0.007972098912524679},
  @{A="sCs"; B="Anxa1"; C="Fpr2"; D="M1"; E=3; F=1; G=77.646779; H=232.940337; I=0.1315130278311194; J=0.1315130278311194; K=3; L=1; M=4.402069333333333; N=13.206208; O=0.2405808023879173; P=0.2405808023879173; Q=341.8065046680106; R=3076.258542012096; S=0.0316395097600752; T=0.03163950976007521},
  @{A="sCs"; B="Anxa1"; C="Fpr2"; D="M2"; E=3; F=1; G=77.646779; H=232.940337; I=0.1315130278311194; J=0.1315130278311194; K=3; L=1; M=12.78643133333333; N=38.359294; O=0.6988008767962779; P=0.6988008767962779; Q=992.8252079380086; R=8935.426871442078; S=0.09190141915851953; T=0.09190141915851956}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r = $r + 1
}

Write-Host "done writing $($rows.Count) rows"